$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add the new "swimlanes" sheet ---
$data = $wb.Sheets.Item(1)
$data.Name = "data"

$swimlanes = $wb.Sheets.Add([System.Reflection.Missing]::Value, $data)
$swimlanes.Name = "swimlanes"

# --- Populate "swimlanes" with the header + artist rows from "data" (B4:E34 -> B3:E33) ---
$src = $data.Range("B4:E34")
$src.Copy()
$swimlanes.Range("B3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Header row (Naam/van/tot/url) rendered bold, same as on "data"
$swimlanes.Range("B3:E3").Font.Bold = $true

# Agostino Chigi's url cell keeps its explicit (non-hyperlink) black font colour
$swimlanes.Range("E15").Font.Color = 0

# Recreate the two hyperlinks (Albrecht Dürer, Niccolò Machiavelli) shifted up one row
$swimlanes.Hyperlinks.Add($swimlanes.Range("E18"), "https://it.wikipedia.org/wiki/Albrecht_D" + [char]0x00FC + "rer")
$swimlanes.Range("E18").Style = "Lien hypertexte"

$swimlanes.Hyperlinks.Add($swimlanes.Range("E23"), "https://it.wikipedia.org/wiki/Niccol" + [char]0x00F2 + "_Machiavelli")
$swimlanes.Range("E23").Style = "Lien hypertexte"

# --- View state: "data" is no longer the selected tab, selection becomes B4:E34 ---
[void]$data.Activate()
$excel.ActiveWindow.Zoom = 157
[void]$data.Range("B4:E34").Select()

# --- View state: "swimlanes" becomes the active/selected tab ---
[void]$swimlanes.Activate()
$excel.ActiveWindow.Zoom = 138
[void]$swimlanes.Range("H9").Select()
